$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J (copy H1 formatting, then set values)
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-19
$iValues = @(9, 9, 9, 9, 9, 8, 9, 8, 9, 8, 9, 8, 9, 6, 7, 7, 9, 8)
$jValues = @(9, 9, 9, 9, 9, 9, 9, 8, 9, 8, 9, 8, 9, 6, 7, 7, 9, 8)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
